$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text while we write the new values,
# then restore the default (Normal) style so no stray number formatting
# is left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "53.979.56"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "2.240.04"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "494.70"
$ws.Range("E5").Value = "  +2.41%  "

$ws.Range("D6").Value = "127.35"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("D9").Value = "2.277.25"
$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("D10").Value = "0.0949"
$ws.Range("E10").Value = "  +3.60%  "

$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  +3.32%  "

$ws.Range("D13").Value = "4.61"
$ws.Range("E13").Value = "  -2.11%  "

$ws.Range("D14").Value = "2.677.91"
$ws.Range("E14").Value = "  +1.85%  "

$ws.Range("D15").Value = "21.77"
$ws.Range("E15").Value = "  +3.52%  "

$ws.Range("D16").Value = "53.992.72"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").Value = "2.289.80"
$ws.Range("E18").Value = "  +2.34%  "

$ws.Range("D19").Value = "10.03"
$ws.Range("E19").Value = "  +4.94%  "

$ws.Range("D20").Value = "4.09"
$ws.Range("E20").Value = "  +3.68%  "

$ws.Range("E21").Value = "  +5.71%  "

$ws.Range("D22").Value = "300.79"
$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").Value = "5.39"
$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D25").Value = "62.52"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("E27").Value = "  +2.49%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.364.24"
$ws.Range("E28").Value = "  +1.10%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.148"
$ws.Range("E29").Value = "  +4.17%  "

$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("D31").Value = "168.66"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").Value = "0.0₃0686"
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("E34").Value = "  +2.48%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("E37").Value = "  +2.27%  "

$ws.Range("D38").Value = "17.59"
$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("D40").Value = "0.863"
$ws.Range("E40").Value = "  +3.40%  "

$ws.Range("E41").Value = "  +4.05%  "

$ws.Range("D42").Value = "35.39"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("E43").Value = "  +3.43%  "

$ws.Range("E44").Value = "  +2.39%  "

$ws.Range("E45").Value = "  +2.14%  "

$ws.Range("D46").Value = "127.95"
$ws.Range("E46").Value = "  +4.66%  "

$ws.Range("D47").Value = "4.76"
$ws.Range("E47").Value = "  +3.07%  "

$ws.Range("D48").Value = "0.0887"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").Value = "0.540"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").Value = "238.08"
$ws.Range("E50").Value = "  +3.10%  "

$ws.Range("E51").Value = "  +2.97%  "

# Restore default styling on the Price column (matches original workbook).
$priceRange.Style = "Normal"
